$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The columns "codeforiati:group-name" (C) and "codeforiati:group-code" (D)
# had their shared-string entries reordered (code now precedes name), which
# swaps the values Excel shows in columns C and D for the header and every
# data row. Reproduce that by swapping the whole C/D column contents.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$colC = $ws.Range("C1:C$lastRow")
$colD = $ws.Range("D1:D$lastRow")

$valuesC = $colC.Value2
$valuesD = $colD.Value2

$colC.Value2 = $valuesD
$colD.Value2 = $valuesC
